$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.919.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.633.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.44"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.258"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.864.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.634.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.564"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.910.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0720"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("E23").Value = "  -3.43%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("E26").Value = "  -1.11%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0481"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.394.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("E35").Value = "  +0.76%  "
$ws.Range("E36").Value = "  +10.21%  "
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("E38").Value = "  +1.20%  "
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.850"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.81%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.79%  "
$ws.Range("E45").Value = "  -1.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.774.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("E47").Value = "  -2.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("E49").Value = "  +1.60%  "
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.99%  "
